$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2 value (minor precision change)
$ws.Range("B2").Value = 0.466746127399636

# Add new rows of data
$ws.Range("A3").Value = 20
$ws.Range("B3").Value = 0.2679005757569375

$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 0.2380682266098698

$ws.Range("A5").Value = 40
$ws.Range("B5").Value = 0.2498407623158637
